# Lab 1 answer sheet - fill in the blanks (Tom Ekshtein's submission)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Name: __________  ->  Name: _Tom Ekshtein_______
#    "Tom Ekshtein" keeps the run's existing Bold formatting and gets a
#    single underline.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("__________", $true, $false, $false, $false, $false, $true, 1, $false, "_Tom Ekshtein_______", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Tom Ekshtein", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.Underline = 1

# ---------------------------------------------------------------------
# 2) " is bit number ____"  ->  " is bit number _31___"
#    "31" is the underlined answer (MSB of a 32-bit dword is bit 31).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" is bit number ____", $true, $false, $false, $false, $false, $true, 1, $false, " is bit number _31___", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("is bit number _31", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.End - 2
$endp = $rng.End
$answerRng = $d.Range($start, $endp)
$answerRng.Font.Underline = 1

# ---------------------------------------------------------------------
# 3) "Is this color value signed or unsigned data?   _________"
#    ->  "...  __no___?____"   ("no" plain, not underlined)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Is this color value signed or unsigned data?   _________", $true, $false, $false, $false, $false, $true, 1, $false, "Is this color value signed or unsigned data?   __no___?____", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "How many bits is needed to represent a data value between 0 to 255?  _______"
#    ->  "...  ___256__?__"   ("256" plain, not underlined)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("How many bits is needed to represent a data value between 0 to 255?  _______", $true, $false, $false, $false, $false, $true, 1, $false, "How many bits is needed to represent a data value between 0 to 255?  ___256__?__", 2) | Out-Null
